$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell, matching the style of the existing header row (E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Populate the new time_taken column for each data row
$ws.Range("F2").Value = "2021-10-05 13:40:47.405968"
$ws.Range("F3").Value = "2021-10-05 13:40:47.405978"
$ws.Range("F4").Value = "2021-10-05 13:40:47.405981"
$ws.Range("F5").Value = "2021-10-05 13:40:47.405984"
$ws.Range("F6").Value = "2021-10-05 13:40:47.405987"
